$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "27.212.37", "  -1.63%  ")
    ,@("Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "1.820.42", "  -2.08%  ")
    ,@("TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "'1.006", "  -1.34%  ")
    ,@("BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "'313.78", "  -2.05%  ")
    ,@("USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "'1.004", "  -1.34%  ")
    ,@("XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "'0.4262", "  -2.31%  ")
    ,@("Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "'0.3664", "  -3.39%  ")
    ,@("OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "'46.01", "  -1.57%  ")
    ,@("Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "'0.07209", "  -2.86%  ")
    ,@("Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "'0.8598", "  -2.76%  ")
    ,@("Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "'20.95", "  -2.93%  ")
    ,@("WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "1.880.01", "  +0.92%  ")
    ,@("Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "'6.646", "  -1.36%  ")
    ,@("TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "'0.07098", "  +0.04%  ")
    ,@("Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "'5.298", "  -3.42%  ")
    ,@("Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "'87.74", "  +0.94%  ")
    ,@("BinanceUSD", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd", "'1.006", "  -1.66%  ")
    ,@("ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "'0.000008855", "  -2.28%  ")
    ,@("Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "'1.005", "  -1.38%  ")
    ,@("Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "'15.01", "  -2.85%  ")
    ,@("WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "27.236.47", "  -1.59%  ")
    ,@("Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "'5.120", "  -3.05%  ")
    ,@("Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "'10.87", "  -2.58%  ")
    ,@("WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "2.062.68", "  -2.02%  ")
    ,@("Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "'2.002", "  -1.63%  ")
    ,@("Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "'152.92", "  -2.79%  ")
    ,@("EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "'18.27", "  -2.37%  ")
    ,@("LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "'2.105", "  +5.40%  ")
    ,@("InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "'5.220", "  -2.67%  ")
    ,@("BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "'115.96", "  -3.74%  ")
    ,@("Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "'0.08875", "  -1.99%  ")
    ,@("ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "'1.191", "  -2.21%  ")
    ,@("ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "'0.7592", "  -1.26%  ")
    ,@("Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "'4.447", "  -2.43%  ")
    ,@("HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "'2.790", "  -8.00%  ")
    ,@("Frax", "https://coinranking.com/coin/KfWtaeV1W+frax-frax", "'1.004", "  -1.46%  ")
    ,@("TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "'1.114", "  -2.35%  ")
    ,@("VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "'0.01953", "  -1.18%  ")
    ,@("Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "'0.05242", "  -1.03%  ")
    ,@("MXToken", "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx", "'2.889", "  +0.79%  ")
    ,@("FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "'7.042", "  +1.31%  ")
    ,@("Algorand", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo", "'0.1675", "  -0.36%  ")
    ,@("TheSandbox", "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand", "'0.5008", "  -3.55%  ")
    ,@("Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "'8.610", "  -0.99%  ")
    ,@("EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "'10.60", "  -1.54%  ")
    ,@("Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "'106.26", "  -3.47%  ")
    ,@("Decentraland", "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana", "'0.4679", "  -1.03%  ")
    ,@("PaxDollar", "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp", "'1.004", "  -1.50%  ")
    ,@("Cronos", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro", "'0.06393", "  -1.55%  ")
    ,@("NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "'1.656", "  -3.31%  ")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 2).Value = $data[$i][0]
    $ws.Cells.Item($r, 3).Value = $data[$i][1]
    $ws.Cells.Item($r, 4).Value = $data[$i][2]
    $ws.Cells.Item($r, 5).Value = $data[$i][3]
}